{"js": "// Change 1: In the Q5 solution paragraph, \"A Q-Q Plot will help us know if the\n// data are normally distributed.\" becomes \"A histogram will help us know if\n// the data are normally distributed.\"\n// Change 2: In the Q10 solution paragraph (\"b. The sample size for\n// non-certified cars...\"), the whole explanatory sentence about the\n// assumption check is rewritten to talk about the sampling distribution of\n// the sample mean and to reference a histogram instead of a Q-Q Plot.\n\nconst body = context.document.body;\n\n// --- Change 1 -----------------------------------------------------------\nconst oldSentence1 =\n  \"10 . A Q-Q Plot will help us know if the data are normally distributed.\";\nconst newSentence1 =\n  \"10 . A histogram will help us know if the data are normally distributed.\";\n\nconst search1 = body.search(oldSentence1, { matchCase: true });\nsearch1.load(\"text\");\nawait context.sync();\n\nif (search1.items.length === 0) {\n  throw new Error(\"Could not find the Q5 Q-Q Plot sentence to update.\");\n}\n\nsearch1.items[0].insertText(newSentence1, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Change 2 -------------------------------------------------------------\nconst oldSentence2 =\n  \"b. The sample size for non-certified cars is large, so we can assume that \" +\n  \"it is normally distributed. With a sample size of n = 24 for certified \" +\n  \"cars, it is likely large enough to assume a normal distribution, but a \" +\n  \"Q-Q Plot will help us know more confidently.\";\nconst newSentence2 =\n  \"b. The sample size for non-certified cars is large, so we can assume that \" +\n  \"the sampling distribution of the sample mean is normally distributed. \" +\n  \"With a sample size of n = 24 for certified cars, it is likely large enough \" +\n  \"to assume a normal distribution for the mean, but a histogram of the data \" +\n  \"will help us know if the data is normally distributed. If it is, then the \" +\n  \"distribution of the sample mean will also be normally distributed.\";\n\nconst search2 = body.search(oldSentence2, { matchCase: true });\nsearch2.load(\"text\");\nawait context.sync();\n\nif (search2.items.length === 0) {\n  throw new Error(\"Could not find the Q10 non-certified cars sentence to update.\");\n}\n\nsearch2.items[0].insertText(newSentence2, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Ryan Cromar email 11/2/2021 1:24 PM\n#\n# Change 1: In the Q5 solution paragraph, \"A Q-Q Plot will help us know if\n# the data are normally distributed.\" becomes \"A histogram will help us\n# know if the data are normally distributed.\"\n#\n# Change 2: In the Q10 solution paragraph (\"b. The sample size for\n# non-certified cars...\"), the explanatory sentence is rewritten to talk\n# about the sampling distribution of the sample mean and to reference a\n# histogram instead of a Q-Q Plot.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($FindText, $ReplaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $FindText\n    $rng.Find.Replacement.Text = $ReplaceText\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n\n    # wdFindContinue = 1, wdReplaceOne = 2\n    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $FindText\"\n    }\n}\n\nReplace-ExactText \"10 . A Q-Q Plot will help us know if the data are normally distributed.\" \"10 . A histogram will help us know if the data are normally distributed.\"\n\n$oldSentence2 = \"b. The sample size for non-certified cars is large, so we can assume that it is normally distributed. With a sample size of n = 24 for certified cars, it is likely large enough to assume a normal distribution, but a Q-Q Plot will help us know more confidently.\"\n$newSentence2 = \"b. The sample size for non-certified cars is large, so we can assume that the sampling distribution of the sample mean is normally distributed. With a sample size of n = 24 for certified cars, it is likely large enough to assume a normal distribution for the mean, but a histogram of the data will help us know if the data is normally distributed. If it is, then the distribution of the sample mean will also be normally distributed.\"\n\nReplace-ExactText $oldSentence2 $newSentence2\n"}
